$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: issue number and reporting week dates ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Crime Complaints data table updates (rows 15-31) ---
$ws.Range("G15").Value = 1
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 2
$ws.Range("C33").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("C33").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("C33").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 33
$ws.Range("K16").Value = 3.125
$ws.Range("L16").Value = 6.451612903225
$ws.Range("M16").Value = 37.5
$ws.Range("N16").Value = -85.897435897435
$ws.Range("C33").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 3
$ws.Range("K33").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 20
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 16.666666666666
$ws.Range("N17").Value = -27.941176470588
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 39
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -26.415094339622
$ws.Range("L18").Value = -49.350649350649
$ws.Range("M18").Value = -20.408163265306
$ws.Range("N18").Value = -91.409691629955
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -29.411764705882
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -32.786885245901
$ws.Range("I19").Value = 233
$ws.Range("J19").Value = 297
$ws.Range("K19").Value = -21.548821548821
$ws.Range("L19").Value = -7.905138339920
$ws.Range("M19").Value = -16.487455197132
$ws.Range("N19").Value = -71.654501216545
$ws.Range("J33").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("J33").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = -80
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -57.142857142857
$ws.Range("L20").Value = -40
$ws.Range("M20").Value = -10
$ws.Range("N20").Value = -96.629213483146
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -40
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -32.142857142857
$ws.Range("I21").Value = 368
$ws.Range("J21").Value = 447
$ws.Range("K21").Value = -17.673378076062
$ws.Range("L21").Value = -13.207547169811
$ws.Range("M21").Value = -6.122448979591
$ws.Range("N21").Value = -80.108108108108
$ws.Range("J33").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("K33").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 50
$ws.Range("M22").Value = 23.529411764705
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = -2.597402597402
$ws.Range("I24").Value = 386
$ws.Range("J24").Value = 380
$ws.Range("K24").Value = 1.578947368421
$ws.Range("L24").Value = -21.384928716904
$ws.Range("M24").Value = 51.968503937007
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 316
$ws.Range("J25").Value = 307
$ws.Range("K25").Value = 2.931596091205
$ws.Range("L25").Value = -23.300970873786
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -39.285714285714
$ws.Range("I26").Value = 84
$ws.Range("J26").Value = 107
$ws.Range("K26").Value = -21.495327102803
$ws.Range("L26").Value = -11.578947368421
$ws.Range("M26").Value = -10.638297872340
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -8
$ws.Range("L28").Value = 9.523809523809
$ws.Range("J33").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = 1
$ws.Range("F31").Value = 4
$ws.Range("C33").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("C33").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("C33").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("I31").Value = 8
$ws.Range("K31").Value = 166.666666666667
$ws.Range("L31").Value = 60

$excel.CutCopyMode = $false
